$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure cells keep their original text format (values look numeric,
# e.g. "28.254.70", "11.00", "0.07410") instead of being auto-converted to numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.254.70"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +3.31%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.924.06"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +3.13%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.006"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -1.50%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.56"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.75%  "
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -1.42%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4855"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.95%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3854"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +3.22%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07410"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.05%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9471"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.93%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.70%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07796"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.19%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.946.38"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +4.35%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.543"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.88%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.676"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.78%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "92.07"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.93%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.006"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.59%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008881"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.93%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -1.33%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "28.258.72"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +3.16%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.04"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +1.52%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.175"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.77%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.158.88"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +2.52%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.00"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +2.53%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -1.67%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "156.38"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.35%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.67"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.118"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +5.30%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "117.23"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.95%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.013"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.31%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08909"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.25%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.30%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.255"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +4.83%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +4.30%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +2.57%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.773"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +2.97%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02058"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.20%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.132"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.49%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05375"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +1.37%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5595"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +4.22%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.045"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +1.29%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "7.106"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.19%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.592"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +2.23%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1539"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.79%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4925"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +1.79%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.74"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.62%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "107.18"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +3.64%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -1.49%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.678"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.29%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "69.64"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +4.07%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06155"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.90%  "
